# Update for March 21
# Adds the March 21, 2020 data row (377 Ontario cases, 23384 tests conducted)
# to the OntarioCoronavirus.csv tracking sheet, together with its source
# hyperlink, and shifts the chart placed below the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 27
$prevRow = $newRow - 1

# --- Data row -------------------------------------------------------------
# A: Date (Excel serial 43911 == 2020-03-21)
$ws.Cells.Item($newRow, 1).Value = 43911
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

# B: Ontario Cases
$ws.Cells.Item($newRow, 2).Value = 377

# C: Tests Conducted
$ws.Cells.Item($newRow, 3).Value = 23384

# D: Days since start -- fill the shared formula down from the row above
$ws.Cells.Item($newRow, 4).Formula = "=A$newRow-`$A`$2"
$ws.Cells.Item($newRow, 4).Style = "Normal"

# E: Day Gap
$ws.Cells.Item($newRow, 5).Formula = "=D$newRow-D$prevRow"
$ws.Cells.Item($newRow, 5).Style = "Normal"

# F: Adjusted Growth
$ws.Cells.Item($newRow, 6).Formula = "=(B$newRow/B$prevRow)^(1/E$newRow)-1"
$ws.Cells.Item($newRow, 6).NumberFormat = $ws.Cells.Item($prevRow, 6).NumberFormat

# G: Source (text + hyperlink, matching the style used by the rest of the column)
$sourceUrl = "https://www.cbc.ca/news/canada/toronto/ontario-hospital-association-pleads-for-social-distancing-1.5505827"
$ws.Cells.Item($newRow, 7).Value = $sourceUrl
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 7), $sourceUrl) | Out-Null
$ws.Cells.Item($prevRow, 7).Copy()
$ws.Cells.Item($newRow, 7).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Move the chart down by one row ---------------------------------------
# The chart sits just below the table; since the table grew by one row, push
# the chart down by the height of the newly-added row so it still starts
# right under the data.
$co = $ws.ChartObjects(1)
$co.Top = $co.Top + $ws.Rows.Item($newRow).Height

$wb.Save()
